$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lung")

# Update the header label and the parameter/variation values on the Lung sheet
$ws.Range("A1").Value = "Day 1, n = 1000"

$ws.Range("B2").Value = 78.85
$ws.Range("C2").Value = 16.6

$ws.Range("B3").Value = 79.79
$ws.Range("C3").Value = 36.45
